$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 800  # H32
$ws.Cells.Item(32, 9).Value = 0  # I32
$ws.Cells.Item(32, 10).Value = 800  # J32
$ws.Cells.Item(32, 11).Value = 0  # K32
$ws.Cells.Item(32, 12).Value = 800  # L32
$ws.Cells.Item(32, 13).Value = ""  # M32
$ws.Cells.Item(32, 14).Value = -1452  # N32
$ws.Cells.Item(94, 8).Value = 5269.2856  # H94
$ws.Cells.Item(94, 9).Value = 5949.1665  # I94
$ws.Cells.Item(94, 11).Value = 5949.1665  # K94
$ws.Cells.Item(94, 13).Value = -5498.1665  # M94
$ws.Cells.Item(100, 8).Value = 1885.8889  # H100
$ws.Cells.Item(100, 9).Value = 877.36365  # I100
$ws.Cells.Item(100, 11).Value = 877.36365  # K100
$ws.Cells.Item(100, 13).Value = -336.36365  # M100
$ws.Cells.Item(112, 8).Value = 1250  # H112
$ws.Cells.Item(112, 9).Value = 1250  # I112
$ws.Cells.Item(112, 10).Value = 0  # J112
$ws.Cells.Item(112, 11).Value = 3750  # K112
$ws.Cells.Item(112, 12).Value = 0  # L112
$ws.Cells.Item(112, 13).Value = -2642  # M112
$ws.Cells.Item(112, 14).Value = ""  # N112
$ws.Cells.Item(113, 8).Value = 2098.2  # H113
$ws.Cells.Item(113, 9).Value = 1997  # I113
$ws.Cells.Item(113, 11).Value = 1997  # K113
$ws.Cells.Item(113, 13).Value = 1257  # M113
$ws.Cells.Item(116, 8).Value = 7915  # H116
$ws.Cells.Item(116, 9).Value = 13137.5  # I116
$ws.Cells.Item(116, 10).Value = 4433.3335  # J116
$ws.Cells.Item(116, 11).Value = 13137.5  # K116
$ws.Cells.Item(116, 12).Value = 4433.3335  # L116
$ws.Cells.Item(116, 13).Value = -9695.5  # M116
$ws.Cells.Item(116, 14).Value = -11317.3335  # N116
$ws.Cells.Item(132, 8).Value = 15104.777  # H132
$ws.Cells.Item(132, 9).Value = 14459.066  # I132
$ws.Cells.Item(132, 10).Value = 18333.334  # J132
$ws.Cells.Item(132, 11).Value = 43377.198  # K132
$ws.Cells.Item(132, 12).Value = 55000.00199999999  # L132
$ws.Cells.Item(132, 13).Value = -40847.198  # M132
$ws.Cells.Item(132, 14).Value = -60060.00199999999  # N132

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7083.5  # H32
$ws.Cells.Item(32, 9).Value = 4616.1304  # I32
$ws.Cells.Item(32, 11).Value = 4616.1304  # K32
$ws.Cells.Item(32, 13).Value = -4329.1304  # M32
$ws.Cells.Item(122, 8).Value = 2129.3635  # H122
$ws.Cells.Item(122, 9).Value = 2192.3  # I122
$ws.Cells.Item(122, 10).Value = 1500  # J122
$ws.Cells.Item(122, 11).Value = 6576.900000000001  # K122
$ws.Cells.Item(122, 12).Value = 4500  # L122
$ws.Cells.Item(122, 13).Value = -4126.900000000001  # M122
$ws.Cells.Item(122, 14).Value = -9400  # N122

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5, 8).Value = 371.5  # H5
$ws.Cells.Item(5, 9).Value = 193  # I5
$ws.Cells.Item(5, 11).Value = 193  # K5
$ws.Cells.Item(5, 13).Value = -80  # M5
$ws.Cells.Item(8, 8).Value = 107.5  # H8
$ws.Cells.Item(8, 9).Value = 140  # I8
$ws.Cells.Item(8, 10).Value = 75  # J8
$ws.Cells.Item(8, 11).Value = 140  # K8
$ws.Cells.Item(8, 12).Value = 75  # L8
$ws.Cells.Item(8, 13).Value = 0  # M8
$ws.Cells.Item(8, 14).Value = -355  # N8
$ws.Cells.Item(10, 8).Value = 241.66667  # H10
$ws.Cells.Item(10, 9).Value = 241.66667  # I10
$ws.Cells.Item(10, 10).Value = 0  # J10
$ws.Cells.Item(10, 11).Value = 241.66667  # K10
$ws.Cells.Item(10, 12).Value = 0  # L10
$ws.Cells.Item(10, 13).Value = -101.66667  # M10
$ws.Cells.Item(10, 14).Value = ""  # N10
$ws.Cells.Item(11, 8).Value = 875  # H11
$ws.Cells.Item(11, 10).Value = 325  # J11
$ws.Cells.Item(11, 12).Value = 325  # L11
$ws.Cells.Item(11, 14).Value = -605  # N11
$ws.Cells.Item(12, 8).Value = 179  # H12
$ws.Cells.Item(12, 9).Value = 221.25  # I12
$ws.Cells.Item(12, 10).Value = 10  # J12
$ws.Cells.Item(12, 11).Value = 221.25  # K12
$ws.Cells.Item(12, 12).Value = 10  # L12
$ws.Cells.Item(12, 13).Value = -53.25  # M12
$ws.Cells.Item(12, 14).Value = -346  # N12
$ws.Cells.Item(14, 8).Value = 862  # H14
$ws.Cells.Item(14, 10).Value = 882.6667  # J14
$ws.Cells.Item(14, 12).Value = 882.6667  # L14
$ws.Cells.Item(14, 14).Value = -1226.6667  # N14
$ws.Cells.Item(23, 8).Value = 15000  # H23
$ws.Cells.Item(23, 9).Value = 10000  # I23
$ws.Cells.Item(23, 11).Value = 10000  # K23
$ws.Cells.Item(23, 13).Value = -9717  # M23
$ws.Cells.Item(86, 8).Value = 3307.8262  # H86
$ws.Cells.Item(86, 9).Value = 2192.5881  # I86
$ws.Cells.Item(86, 10).Value = 6467.6665  # J86
$ws.Cells.Item(86, 11).Value = 2192.5881  # K86
$ws.Cells.Item(86, 12).Value = 6467.6665  # L86
$ws.Cells.Item(86, 13).Value = -1069.5881  # M86
$ws.Cells.Item(86, 14).Value = -8713.666499999999  # N86
$ws.Cells.Item(89, 8).Value = 3307.8262  # H89
$ws.Cells.Item(89, 9).Value = 2192.5881  # I89
$ws.Cells.Item(89, 10).Value = 6467.6665  # J89
$ws.Cells.Item(89, 11).Value = 10962.9405  # K89
$ws.Cells.Item(89, 12).Value = 32338.3325  # L89
$ws.Cells.Item(89, 13).Value = -5346.940500000001  # M89
$ws.Cells.Item(89, 14).Value = -43570.3325  # N89

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 777.5  # H2
$ws.Cells.Item(2, 10).Value = 897.5  # J2
$ws.Cells.Item(2, 12).Value = 897.5  # L2
$ws.Cells.Item(2, 14).Value = -1123.5  # N2
$ws.Cells.Item(5, 8).Value = 608.3  # H5
$ws.Cells.Item(5, 9).Value = 229.125  # I5
$ws.Cells.Item(5, 10).Value = 2125  # J5
$ws.Cells.Item(5, 11).Value = 229.125  # K5
$ws.Cells.Item(5, 12).Value = 2125  # L5
$ws.Cells.Item(5, 13).Value = -117.125  # M5
$ws.Cells.Item(5, 14).Value = -2349  # N5
$ws.Cells.Item(11, 8).Value = 453  # H11
$ws.Cells.Item(11, 10).Value = 21.333334  # J11
$ws.Cells.Item(11, 12).Value = 21.333334  # L11
$ws.Cells.Item(11, 14).Value = -301.333334  # N11
$ws.Cells.Item(62, 8).Value = 4999  # H62
$ws.Cells.Item(62, 10).Value = 4999  # J62
$ws.Cells.Item(62, 12).Value = 4999  # L62
$ws.Cells.Item(62, 14).Value = -6247  # N62
$ws.Cells.Item(65, 8).Value = 4999  # H65
$ws.Cells.Item(65, 10).Value = 4999  # J65
$ws.Cells.Item(65, 12).Value = 24995  # L65
$ws.Cells.Item(65, 14).Value = -31235  # N65
$ws.Cells.Item(99, 8).Value = 2612.818  # H99
$ws.Cells.Item(99, 9).Value = 2092.5  # I99
$ws.Cells.Item(99, 11).Value = 2092.5  # K99
$ws.Cells.Item(99, 13).Value = -594.5  # M99
$ws.Cells.Item(105, 8).Value = 1345.2858  # H105
$ws.Cells.Item(105, 9).Value = 1411.1666  # I105
$ws.Cells.Item(105, 10).Value = 950  # J105
$ws.Cells.Item(105, 11).Value = 1411.1666  # K105
$ws.Cells.Item(105, 12).Value = 950  # L105
$ws.Cells.Item(105, 13).Value = 335.8334  # M105
$ws.Cells.Item(105, 14).Value = -4444  # N105
$ws.Cells.Item(126, 8).Value = 2612.818  # H126
$ws.Cells.Item(126, 9).Value = 2092.5  # I126
$ws.Cells.Item(126, 11).Value = 6277.5  # K126
$ws.Cells.Item(126, 13).Value = -3807.5  # M126

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 83.59999999999999  # H2
$ws.Cells.Item(2, 10).Value = 39.333332  # J2
$ws.Cells.Item(2, 12).Value = 235.999992  # L2
$ws.Cells.Item(2, 14).Value = -461.999992  # N2
$ws.Cells.Item(138, 8).Value = 6922.273  # H138
$ws.Cells.Item(138, 9).Value = 2211.25  # I138
$ws.Cells.Item(138, 10).Value = 9614.286  # J138
$ws.Cells.Item(138, 11).Value = 6633.75  # K138
$ws.Cells.Item(138, 12).Value = 28842.858  # L138
$ws.Cells.Item(138, 13).Value = -1493.75  # M138
$ws.Cells.Item(138, 14).Value = -39122.858  # N138

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(6, 8).Value = 116.333336  # H6
$ws.Cells.Item(6, 9).Value = 99  # I6
$ws.Cells.Item(6, 10).Value = 125  # J6
$ws.Cells.Item(6, 11).Value = 99  # K6
$ws.Cells.Item(6, 12).Value = 125  # L6
$ws.Cells.Item(6, 13).Value = 14  # M6
$ws.Cells.Item(6, 14).Value = -351  # N6
$ws.Cells.Item(7, 8).Value = 21546046  # H7
$ws.Cells.Item(7, 10).Value = 40000830  # J7
$ws.Cells.Item(7, 12).Value = 40000830  # L7
$ws.Cells.Item(7, 14).Value = -40001054  # N7
$ws.Cells.Item(8, 8).Value = 21546046  # H8
$ws.Cells.Item(8, 10).Value = 40000830  # J8
$ws.Cells.Item(8, 12).Value = 40000830  # L8
$ws.Cells.Item(8, 14).Value = -40001108  # N8
$ws.Cells.Item(9, 8).Value = 218.88889  # H9
$ws.Cells.Item(9, 9).Value = 252.85715  # I9
$ws.Cells.Item(9, 11).Value = 252.85715  # K9
$ws.Cells.Item(9, 13).Value = -82.85714999999999  # M9
$ws.Cells.Item(13, 8).Value = 200.05882  # H13
$ws.Cells.Item(13, 10).Value = 208.13333  # J13
$ws.Cells.Item(13, 12).Value = 208.13333  # L13
$ws.Cells.Item(13, 14).Value = -486.13333  # N13
$ws.Cells.Item(16, 8).Value = 116.333336  # H16
$ws.Cells.Item(16, 9).Value = 99  # I16
$ws.Cells.Item(16, 10).Value = 125  # J16
$ws.Cells.Item(16, 11).Value = 99  # K16
$ws.Cells.Item(16, 12).Value = 125  # L16
$ws.Cells.Item(16, 13).Value = 151  # M16
$ws.Cells.Item(16, 14).Value = -625  # N16
$ws.Cells.Item(17, 8).Value = 25204.5  # H17
$ws.Cells.Item(17, 10).Value = 25204.5  # J17
$ws.Cells.Item(17, 12).Value = 25204.5  # L17
$ws.Cells.Item(17, 14).Value = -25540.5  # N17
$ws.Cells.Item(19, 8).Value = 500  # H19
$ws.Cells.Item(19, 10).Value = 500  # J19
$ws.Cells.Item(19, 12).Value = 500  # L19
$ws.Cells.Item(19, 14).Value = -1076  # N19
$ws.Cells.Item(22, 8).Value = 504.5  # H22
$ws.Cells.Item(22, 9).Value = 504.5  # I22
$ws.Cells.Item(22, 10).Value = 0  # J22
$ws.Cells.Item(22, 11).Value = 504.5  # K22
$ws.Cells.Item(22, 12).Value = 0  # L22
$ws.Cells.Item(22, 13).Value = 24.5  # M22
$ws.Cells.Item(22, 14).Value = ""  # N22
$ws.Cells.Item(23, 8).Value = 778  # H23
$ws.Cells.Item(23, 10).Value = 778  # J23
$ws.Cells.Item(23, 12).Value = 778  # L23
$ws.Cells.Item(23, 14).Value = -1224  # N23
$ws.Cells.Item(25, 8).Value = 2066.3333  # H25
$ws.Cells.Item(25, 9).Value = 4000  # I25
$ws.Cells.Item(25, 10).Value = 1099.5  # J25
$ws.Cells.Item(25, 11).Value = 4000  # K25
$ws.Cells.Item(25, 12).Value = 1099.5  # L25
$ws.Cells.Item(25, 13).Value = -3471  # M25
$ws.Cells.Item(25, 14).Value = -2157.5  # N25
$ws.Cells.Item(102, 8).Value = 1623  # H102
$ws.Cells.Item(102, 9).Value = 1917.2222  # I102
$ws.Cells.Item(102, 11).Value = 1917.2222  # K102
$ws.Cells.Item(102, 13).Value = -295.2221999999999  # M102
$ws.Cells.Item(126, 8).Value = 6813.5  # H126
$ws.Cells.Item(126, 9).Value = 6582.3335  # I126
$ws.Cells.Item(126, 11).Value = 19747.0005  # K126
$ws.Cells.Item(126, 13).Value = -17277.0005  # M126
$ws.Cells.Item(128, 8).Value = 0  # H128
$ws.Cells.Item(128, 10).Value = 0  # J128
$ws.Cells.Item(128, 12).Value = 0  # L128
$ws.Cells.Item(128, 14).Value = ""  # N128

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8926.286  # H7
$ws.Cells.Item(7, 10).Value = 9199  # J7
$ws.Cells.Item(7, 12).Value = 9199  # L7
$ws.Cells.Item(7, 14).Value = -9423  # N7
$ws.Cells.Item(40, 8).Value = 5009.8667  # H40
$ws.Cells.Item(40, 9).Value = 4970.9287  # I40
$ws.Cells.Item(40, 11).Value = 4970.9287  # K40
$ws.Cells.Item(40, 13).Value = -4834.9287  # M40
$ws.Cells.Item(82, 8).Value = 3660.7856  # H82
$ws.Cells.Item(82, 9).Value = 2269.6667  # I82
$ws.Cells.Item(82, 10).Value = 4704.125  # J82
$ws.Cells.Item(82, 11).Value = 2269.6667  # K82
$ws.Cells.Item(82, 12).Value = 4704.125  # L82
$ws.Cells.Item(82, 13).Value = -1908.6667  # M82
$ws.Cells.Item(82, 14).Value = -5426.125  # N82
$ws.Cells.Item(85, 8).Value = 3660.7856  # H85
$ws.Cells.Item(85, 9).Value = 2269.6667  # I85
$ws.Cells.Item(85, 10).Value = 4704.125  # J85
$ws.Cells.Item(85, 11).Value = 2269.6667  # K85
$ws.Cells.Item(85, 12).Value = 4704.125  # L85
$ws.Cells.Item(85, 13).Value = -1021.6667  # M85
$ws.Cells.Item(85, 14).Value = -7200.125  # N85
$ws.Cells.Item(122, 8).Value = 3112.3333  # H122
$ws.Cells.Item(122, 9).Value = 3002  # I122
$ws.Cells.Item(122, 11).Value = 9006  # K122
$ws.Cells.Item(122, 13).Value = -6556  # M122
$ws.Cells.Item(126, 8).Value = 8926.286  # H126
$ws.Cells.Item(126, 10).Value = 9199  # J126
$ws.Cells.Item(126, 12).Value = 27597  # L126
$ws.Cells.Item(126, 14).Value = -32537  # N126
$ws.Cells.Item(136, 8).Value = 2799.75  # H136
$ws.Cells.Item(136, 9).Value = 2466  # I136
$ws.Cells.Item(136, 11).Value = 7398  # K136
$ws.Cells.Item(136, 13).Value = -4848  # M136

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1506.7142  # H100
$ws.Cells.Item(100, 9).Value = 1506.7142  # I100
$ws.Cells.Item(100, 11).Value = 3013.4284  # K100
$ws.Cells.Item(100, 13).Value = -2472.4284  # M100
$ws.Cells.Item(126, 8).Value = 4017.6316  # H126
$ws.Cells.Item(126, 9).Value = 1673.6  # I126
$ws.Cells.Item(126, 11).Value = 5020.799999999999  # K126
$ws.Cells.Item(126, 13).Value = -2550.799999999999  # M126
